$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-01-25 Saturday" "2025-01-26 Sunday"
Replace-Text "12×47=" "32×81="
Replace-Text "96×93=" "43×47="
Replace-Text "76×29=" "13×11="
Replace-Text "26×83=" "52×12="
Replace-Text "64×53=" "62×85="
Replace-Text "52×78=" "50×83="
Replace-Text "42×21=" "67×58="
Replace-Text "21×38=" "65×69="
Replace-Text "85×60=" "91×33="
Replace-Text "27×91=" "80×88="
Replace-Text "64×73=" "77×39="
Replace-Text "63×95=" "25×14="
Replace-Text "68×25=" "98×42="
Replace-Text "52×37=" "76×32="
Replace-Text "26×11=" "13×66="
Replace-Text "26×52=" "55×83="
Replace-Text "98×85=" "92×87="
Replace-Text "34×40=" "98×50="
Replace-Text "13×23=" "54×63="
Replace-Text "82×94=" "13×21="
Replace-Text "31×21=" "70×91="
Replace-Text "92×79=" "44×71="
Replace-Text "49×40=" "28×16="
Replace-Text "81×80=" "46×76="
Replace-Text "94×22=" "68×30="
